$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

function Set-Cell($ws, $r, $c, $val) {
    $ws.Cells.Item($r, $c).Value = $val
}

# --- games sheet: season 2023, weeks 1-3 results ---
Set-Cell $ws1 285 4 'GB'
Set-Cell $ws1 285 5 'CHI'
Set-Cell $ws1 285 6 38
Set-Cell $ws1 285 7 20
Set-Cell $ws1 285 8 'CHI'

Set-Cell $ws1 286 4 'PHI'
Set-Cell $ws1 286 5 'NE'
Set-Cell $ws1 286 6 25
Set-Cell $ws1 286 7 20
Set-Cell $ws1 286 8 'NE'

Set-Cell $ws1 287 4 'DAL'
Set-Cell $ws1 287 5 'NYG'
Set-Cell $ws1 287 6 40
Set-Cell $ws1 287 7 0
Set-Cell $ws1 287 8 'NYG'

Set-Cell $ws1 288 4 'BUF'
Set-Cell $ws1 288 5 'NYJ'
Set-Cell $ws1 288 6 16
Set-Cell $ws1 288 7 22
Set-Cell $ws1 288 8 'NYJ'

Set-Cell $ws1 289 2 2
$ws1.Cells.Item(280, 3).Copy($ws1.Cells.Item(289, 3))
Set-Cell $ws1 289 3 45183
Set-Cell $ws1 289 4 'MIN'
Set-Cell $ws1 289 5 'PHI'
Set-Cell $ws1 289 6 28
Set-Cell $ws1 289 7 34
Set-Cell $ws1 289 8 'PHI'

Set-Cell $ws1 290 2 2
$ws1.Cells.Item(280, 3).Copy($ws1.Cells.Item(290, 3))
Set-Cell $ws1 290 3 45186
Set-Cell $ws1 290 4 'IND'
Set-Cell $ws1 290 5 'HOU'
Set-Cell $ws1 290 6 31
Set-Cell $ws1 290 7 20
Set-Cell $ws1 290 8 'HOU'

Set-Cell $ws1 291 2 2
$ws1.Cells.Item(280, 3).Copy($ws1.Cells.Item(291, 3))
Set-Cell $ws1 291 3 45186
Set-Cell $ws1 291 4 'LAC'
Set-Cell $ws1 291 5 'TEN'
Set-Cell $ws1 291 6 24
Set-Cell $ws1 291 7 27
Set-Cell $ws1 291 8 'TEN'

Set-Cell $ws1 292 2 2
$ws1.Cells.Item(280, 3).Copy($ws1.Cells.Item(292, 3))
Set-Cell $ws1 292 3 45186
Set-Cell $ws1 292 4 'SEA'
Set-Cell $ws1 292 5 'DET'
Set-Cell $ws1 292 6 37
Set-Cell $ws1 292 7 31
Set-Cell $ws1 292 8 'DET'

Set-Cell $ws1 293 2 2
$ws1.Cells.Item(280, 3).Copy($ws1.Cells.Item(293, 3))
Set-Cell $ws1 293 3 45186
Set-Cell $ws1 293 4 'GB'
Set-Cell $ws1 293 5 'ATL'
Set-Cell $ws1 293 6 24
Set-Cell $ws1 293 7 25
Set-Cell $ws1 293 8 'ATL'

Set-Cell $ws1 294 2 2
$ws1.Cells.Item(280, 3).Copy($ws1.Cells.Item(294, 3))
Set-Cell $ws1 294 3 45186
Set-Cell $ws1 294 4 'LV'
Set-Cell $ws1 294 5 'BUF'
Set-Cell $ws1 294 6 10
Set-Cell $ws1 294 7 38
Set-Cell $ws1 294 8 'BUF'

Set-Cell $ws1 295 2 2
$ws1.Cells.Item(280, 3).Copy($ws1.Cells.Item(295, 3))
Set-Cell $ws1 295 3 45186
Set-Cell $ws1 295 4 'BAL'
Set-Cell $ws1 295 5 'CIN'
Set-Cell $ws1 295 6 27
Set-Cell $ws1 295 7 24
Set-Cell $ws1 295 8 'CIN'

Set-Cell $ws1 296 2 2
$ws1.Cells.Item(280, 3).Copy($ws1.Cells.Item(296, 3))
Set-Cell $ws1 296 3 45186
Set-Cell $ws1 296 4 'KC'
Set-Cell $ws1 296 5 'JAX'
Set-Cell $ws1 296 6 17
Set-Cell $ws1 296 7 9
Set-Cell $ws1 296 8 'JAX'

Set-Cell $ws1 297 2 2
$ws1.Cells.Item(280, 3).Copy($ws1.Cells.Item(297, 3))
Set-Cell $ws1 297 3 45186
Set-Cell $ws1 297 4 'CHI'
Set-Cell $ws1 297 5 'TB'
Set-Cell $ws1 297 6 17
Set-Cell $ws1 297 7 27
$ws1.Cells.Item(297, 8).Formula = "=E297"

Set-Cell $ws1 298 2 2
$ws1.Cells.Item(280, 3).Copy($ws1.Cells.Item(298, 3))
Set-Cell $ws1 298 3 45186
Set-Cell $ws1 298 4 'NYG'
Set-Cell $ws1 298 5 'ARI'
Set-Cell $ws1 298 6 31
Set-Cell $ws1 298 7 28
Set-Cell $ws1 298 8 'ARI'

Set-Cell $ws1 299 2 2
$ws1.Cells.Item(280, 3).Copy($ws1.Cells.Item(299, 3))
Set-Cell $ws1 299 3 45186
Set-Cell $ws1 299 4 'SF'
Set-Cell $ws1 299 5 'LA'
Set-Cell $ws1 299 6 30
Set-Cell $ws1 299 7 23
Set-Cell $ws1 299 8 'LA'

Set-Cell $ws1 300 2 2
$ws1.Cells.Item(280, 3).Copy($ws1.Cells.Item(300, 3))
Set-Cell $ws1 300 3 45186
Set-Cell $ws1 300 4 'WAS'
Set-Cell $ws1 300 5 'DEN'
Set-Cell $ws1 300 6 35
Set-Cell $ws1 300 7 33
Set-Cell $ws1 300 8 'DEN'

Set-Cell $ws1 301 2 2
$ws1.Cells.Item(280, 3).Copy($ws1.Cells.Item(301, 3))
Set-Cell $ws1 301 3 45186
Set-Cell $ws1 301 4 'NYJ'
Set-Cell $ws1 301 5 'DAL'
Set-Cell $ws1 301 6 10
Set-Cell $ws1 301 7 30
Set-Cell $ws1 301 8 'DAL'

Set-Cell $ws1 302 2 2
$ws1.Cells.Item(280, 3).Copy($ws1.Cells.Item(302, 3))
Set-Cell $ws1 302 3 45186
Set-Cell $ws1 302 4 'MIA'
Set-Cell $ws1 302 5 'NE'
Set-Cell $ws1 302 6 24
Set-Cell $ws1 302 7 17
Set-Cell $ws1 302 8 'NE'

Set-Cell $ws1 303 2 2
$ws1.Cells.Item(280, 3).Copy($ws1.Cells.Item(303, 3))
Set-Cell $ws1 303 3 45187
Set-Cell $ws1 303 4 'NO'
Set-Cell $ws1 303 5 'CAR'
Set-Cell $ws1 303 6 20
Set-Cell $ws1 303 7 17
Set-Cell $ws1 303 8 'CAR'

Set-Cell $ws1 304 2 2
$ws1.Cells.Item(280, 3).Copy($ws1.Cells.Item(304, 3))
Set-Cell $ws1 304 3 45187
Set-Cell $ws1 304 4 'CLE'
Set-Cell $ws1 304 5 'PIT'
Set-Cell $ws1 304 6 22
Set-Cell $ws1 304 7 26
Set-Cell $ws1 304 8 'PIT'

Set-Cell $ws1 305 2 3
$ws1.Cells.Item(280, 3).Copy($ws1.Cells.Item(305, 3))
Set-Cell $ws1 305 3 45190
Set-Cell $ws1 305 4 'NYG'
Set-Cell $ws1 305 5 'SF'
Set-Cell $ws1 305 6 12
Set-Cell $ws1 305 7 30
Set-Cell $ws1 305 8 'SF'

Set-Cell $ws1 306 2 3
$ws1.Cells.Item(280, 3).Copy($ws1.Cells.Item(306, 3))
Set-Cell $ws1 306 3 45193
Set-Cell $ws1 306 4 'NO'
Set-Cell $ws1 306 5 'GB'
Set-Cell $ws1 306 6 17
Set-Cell $ws1 306 7 18
Set-Cell $ws1 306 8 'GB'

Set-Cell $ws1 307 2 3
$ws1.Cells.Item(280, 3).Copy($ws1.Cells.Item(307, 3))
Set-Cell $ws1 307 3 45193
Set-Cell $ws1 307 4 'LAC'
Set-Cell $ws1 307 5 'MIN'
Set-Cell $ws1 307 6 28
Set-Cell $ws1 307 7 24
Set-Cell $ws1 307 8 'MIN'

Set-Cell $ws1 308 2 3
$ws1.Cells.Item(280, 3).Copy($ws1.Cells.Item(308, 3))
Set-Cell $ws1 308 3 45193
Set-Cell $ws1 308 4 'NE'
Set-Cell $ws1 308 5 'NYJ'
Set-Cell $ws1 308 6 15
Set-Cell $ws1 308 7 10
Set-Cell $ws1 308 8 'NYJ'

Set-Cell $ws1 309 2 3
$ws1.Cells.Item(280, 3).Copy($ws1.Cells.Item(309, 3))
Set-Cell $ws1 309 3 45193
Set-Cell $ws1 309 4 'DEN'
Set-Cell $ws1 309 5 'MIA'
Set-Cell $ws1 309 6 20
Set-Cell $ws1 309 7 70
Set-Cell $ws1 309 8 'MIA'

Set-Cell $ws1 310 2 3
$ws1.Cells.Item(280, 3).Copy($ws1.Cells.Item(310, 3))
Set-Cell $ws1 310 3 45193
Set-Cell $ws1 310 4 'IND'
Set-Cell $ws1 310 5 'BAL'
Set-Cell $ws1 310 6 22
Set-Cell $ws1 310 7 19
Set-Cell $ws1 310 8 'BAL'

Set-Cell $ws1 311 2 3
$ws1.Cells.Item(280, 3).Copy($ws1.Cells.Item(311, 3))
Set-Cell $ws1 311 3 45193
Set-Cell $ws1 311 4 'ATL'
Set-Cell $ws1 311 5 'DET'
Set-Cell $ws1 311 6 6
Set-Cell $ws1 311 7 20
Set-Cell $ws1 311 8 'DET'

Set-Cell $ws1 312 2 3
$ws1.Cells.Item(280, 3).Copy($ws1.Cells.Item(312, 3))
Set-Cell $ws1 312 3 45193
Set-Cell $ws1 312 4 'BUF'
Set-Cell $ws1 312 5 'WAS'
Set-Cell $ws1 312 6 37
Set-Cell $ws1 312 7 3
Set-Cell $ws1 312 8 'WAS'

Set-Cell $ws1 313 2 3
$ws1.Cells.Item(280, 3).Copy($ws1.Cells.Item(313, 3))
Set-Cell $ws1 313 3 45193
Set-Cell $ws1 313 4 'HOU'
Set-Cell $ws1 313 5 'JAX'
Set-Cell $ws1 313 6 37
Set-Cell $ws1 313 7 17
Set-Cell $ws1 313 8 'JAX'

Set-Cell $ws1 314 2 3
$ws1.Cells.Item(280, 3).Copy($ws1.Cells.Item(314, 3))
Set-Cell $ws1 314 3 45193
Set-Cell $ws1 314 4 'TEN'
Set-Cell $ws1 314 5 'CLE'
Set-Cell $ws1 314 6 3
Set-Cell $ws1 314 7 27
Set-Cell $ws1 314 8 'CLE'

Set-Cell $ws1 315 2 3
$ws1.Cells.Item(280, 3).Copy($ws1.Cells.Item(315, 3))
Set-Cell $ws1 315 3 45193
Set-Cell $ws1 315 4 'CAR'
Set-Cell $ws1 315 5 'SEA'
Set-Cell $ws1 315 6 27
Set-Cell $ws1 315 7 37
Set-Cell $ws1 315 8 'SEA'

Set-Cell $ws1 316 2 3
$ws1.Cells.Item(280, 3).Copy($ws1.Cells.Item(316, 3))
Set-Cell $ws1 316 3 45193
Set-Cell $ws1 316 4 'DAL'
Set-Cell $ws1 316 5 'ARI'
Set-Cell $ws1 316 6 16
Set-Cell $ws1 316 7 28
Set-Cell $ws1 316 8 'ARI'

Set-Cell $ws1 317 2 3
$ws1.Cells.Item(280, 3).Copy($ws1.Cells.Item(317, 3))
Set-Cell $ws1 317 3 45193
Set-Cell $ws1 317 4 'CHI'
Set-Cell $ws1 317 5 'KC'
Set-Cell $ws1 317 6 10
Set-Cell $ws1 317 7 41
Set-Cell $ws1 317 8 'KC'

Set-Cell $ws1 318 2 3
$ws1.Cells.Item(280, 3).Copy($ws1.Cells.Item(318, 3))
Set-Cell $ws1 318 3 45193
Set-Cell $ws1 318 4 'PIT'
Set-Cell $ws1 318 5 'LV'
Set-Cell $ws1 318 6 23
Set-Cell $ws1 318 7 18
Set-Cell $ws1 318 8 'LV'

Set-Cell $ws1 319 2 3
$ws1.Cells.Item(280, 3).Copy($ws1.Cells.Item(319, 3))
Set-Cell $ws1 319 3 45194
Set-Cell $ws1 319 4 'PHI'
Set-Cell $ws1 319 5 'TB'
Set-Cell $ws1 319 6 25
Set-Cell $ws1 319 7 11
Set-Cell $ws1 319 8 'TB'

Set-Cell $ws1 320 2 3
$ws1.Cells.Item(280, 3).Copy($ws1.Cells.Item(320, 3))
Set-Cell $ws1 320 3 45194
Set-Cell $ws1 320 4 'LA'
Set-Cell $ws1 320 5 'CIN'
Set-Cell $ws1 320 6 16
Set-Cell $ws1 320 7 19
Set-Cell $ws1 320 8 'CIN'

# --- games sheet: season 2023, week 4 (schedule only, no results yet) ---
Set-Cell $ws1 321 2 4
Set-Cell $ws1 322 2 4
Set-Cell $ws1 323 2 4
Set-Cell $ws1 324 2 4
Set-Cell $ws1 325 2 4
Set-Cell $ws1 326 2 4
Set-Cell $ws1 327 2 4
Set-Cell $ws1 328 2 4
Set-Cell $ws1 329 2 4
Set-Cell $ws1 330 2 4
Set-Cell $ws1 331 2 4
Set-Cell $ws1 332 2 4
Set-Cell $ws1 333 2 4
Set-Cell $ws1 334 2 4
Set-Cell $ws1 335 2 4
Set-Cell $ws1 336 2 4
# row 321 retains an empty but date-styled C cell
$ws1.Cells.Item(280, 3).Copy($ws1.Cells.Item(321, 3))
$ws1.Cells.Item(321, 3).ClearContents()

# --- 2025 schedule sheet: season 2025 week 1 ---
Set-Cell $ws2 2 2 1
$ws1.Cells.Item(280, 3).Copy($ws2.Cells.Item(2, 3))
Set-Cell $ws2 2 3 45904
Set-Cell $ws2 2 4 'DAL'
Set-Cell $ws2 2 5 'PHI'

Set-Cell $ws2 3 2 1
$ws1.Cells.Item(280, 3).Copy($ws2.Cells.Item(3, 3))
Set-Cell $ws2 3 3 45905
Set-Cell $ws2 3 4 'KC'
Set-Cell $ws2 3 5 'LAC'

Set-Cell $ws2 4 2 1
$ws1.Cells.Item(280, 3).Copy($ws2.Cells.Item(4, 3))
Set-Cell $ws2 4 3 45907
Set-Cell $ws2 4 4 'TB'
Set-Cell $ws2 4 5 'ATL'

Set-Cell $ws2 5 2 1
$ws1.Cells.Item(280, 3).Copy($ws2.Cells.Item(5, 3))
Set-Cell $ws2 5 3 45907
Set-Cell $ws2 5 4 'CIN'
Set-Cell $ws2 5 5 'CLE'

Set-Cell $ws2 6 2 1
$ws1.Cells.Item(280, 3).Copy($ws2.Cells.Item(6, 3))
Set-Cell $ws2 6 3 45907
Set-Cell $ws2 6 4 'MIA'
Set-Cell $ws2 6 5 'IND'

Set-Cell $ws2 7 2 1
$ws1.Cells.Item(280, 3).Copy($ws2.Cells.Item(7, 3))
Set-Cell $ws2 7 3 45907
Set-Cell $ws2 7 4 'CAR'
Set-Cell $ws2 7 5 'JAX'

Set-Cell $ws2 8 2 1
$ws1.Cells.Item(280, 3).Copy($ws2.Cells.Item(8, 3))
Set-Cell $ws2 8 3 45907
Set-Cell $ws2 8 4 'LV'
Set-Cell $ws2 8 5 'NE'

Set-Cell $ws2 9 2 1
$ws1.Cells.Item(280, 3).Copy($ws2.Cells.Item(9, 3))
Set-Cell $ws2 9 3 45907
Set-Cell $ws2 9 4 'ARI'
Set-Cell $ws2 9 5 'NO'

Set-Cell $ws2 10 2 1
$ws1.Cells.Item(280, 3).Copy($ws2.Cells.Item(10, 3))
Set-Cell $ws2 10 3 45907
Set-Cell $ws2 10 4 'PIT'
Set-Cell $ws2 10 5 'NYJ'

Set-Cell $ws2 11 2 1
$ws1.Cells.Item(280, 3).Copy($ws2.Cells.Item(11, 3))
Set-Cell $ws2 11 3 45907
Set-Cell $ws2 11 4 'NYG'
Set-Cell $ws2 11 5 'WAS'

Set-Cell $ws2 12 2 1
$ws1.Cells.Item(280, 3).Copy($ws2.Cells.Item(12, 3))
Set-Cell $ws2 12 3 45907
Set-Cell $ws2 12 4 'TEN'
Set-Cell $ws2 12 5 'DEN'

Set-Cell $ws2 13 2 1
$ws1.Cells.Item(280, 3).Copy($ws2.Cells.Item(13, 3))
Set-Cell $ws2 13 3 45907
Set-Cell $ws2 13 4 'SF'
Set-Cell $ws2 13 5 'SEA'

Set-Cell $ws2 14 2 1
$ws1.Cells.Item(280, 3).Copy($ws2.Cells.Item(14, 3))
Set-Cell $ws2 14 3 45907
Set-Cell $ws2 14 4 'DET'
Set-Cell $ws2 14 5 'GB'

Set-Cell $ws2 15 2 1
$ws1.Cells.Item(280, 3).Copy($ws2.Cells.Item(15, 3))
Set-Cell $ws2 15 3 45907
Set-Cell $ws2 15 4 'HOU'
Set-Cell $ws2 15 5 'LA'

Set-Cell $ws2 16 2 1
$ws1.Cells.Item(280, 3).Copy($ws2.Cells.Item(16, 3))
Set-Cell $ws2 16 3 45907
Set-Cell $ws2 16 4 'BAL'
Set-Cell $ws2 16 5 'BUF'

Set-Cell $ws2 17 2 1
$ws1.Cells.Item(280, 3).Copy($ws2.Cells.Item(17, 3))
Set-Cell $ws2 17 3 45908
Set-Cell $ws2 17 4 'MIN'
Set-Cell $ws2 17 5 'CHI'

# shared formula across F3:F17 (F2 gets its own non-shared formula)
$ws2.Cells.Item(2, 6).Formula = "=E2"
$ws2.Range("F3:F17").Formula = "=E3"

# --- sheet view / selection state ---
$ws2.Activate()
$ws2.Range("C19").Select()
$ws1.Activate()
$ws1.Range("B321:B336").Select()
